# Auto-generated edit script: apply cryptos.xlsx value updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.153.51'
$ws.Range("E2").Value = '''  +1.83%  '

$ws.Range("D3").Value = '''3.906.22'
$ws.Range("E3").Value = '''  +1.20%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '''  +0.14%  '

$ws.Range("D5").Value = '''484.21'
$ws.Range("E5").Value = '''  +3.17%  '

$ws.Range("D6").Value = '''146.32'
$ws.Range("E6").Value = '''  +1.40%  '

$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '''  +2.17%  '

$ws.Range("E8").Value = '''  -0.10%  '

$ws.Range("E9").Value = '''  +3.30%  '

$ws.Range("D10").Value = '''0.174'
$ws.Range("E10").Value = '''  +7.09%  '

$ws.Range("D11").Value = '''0.0000354'
$ws.Range("E11").Value = '''  +4.14%  '

$ws.Range("D12").Value = '''43.06'
$ws.Range("E12").Value = '''  +2.93%  '

$ws.Range("D13").Value = '''10.68'
$ws.Range("E13").Value = '''  +5.61%  '

$ws.Range("D14").Value = '''4.525.22'
$ws.Range("E14").Value = '''  +1.34%  '

$ws.Range("D15").Value = '''3.895.00'
$ws.Range("E15").Value = '''  -0.65%  '

$ws.Range("D16").Value = '''14.37'
$ws.Range("E16").Value = '''  -0.15%  '

$ws.Range("E17").Value = '''  -0.25%  '

$ws.Range("D18").Value = '''20.23'
$ws.Range("E18").Value = '''  +3.58%  '

$ws.Range("E19").Value = '''  +2.37%  '

$ws.Range("D20").Value = '''68.205.79'
$ws.Range("E20").Value = '''  +1.83%  '

$ws.Range("D21").Value = '''429.92'
$ws.Range("E21").Value = '''  +0.68%  '

$ws.Range("D22").Value = '''3.56'
$ws.Range("E22").Value = '''  +8.52%  '

$ws.Range("D23").Value = '''15.00'
$ws.Range("E23").Value = '''  +5.31%  '

$ws.Range("D24").Value = '''89.13'
$ws.Range("E24").Value = '''  +3.94%  '

$ws.Range("D25").Value = '''11.64'
$ws.Range("E25").Value = '''  +14.00%  '

$ws.Range("D26").Value = '''3.71'
$ws.Range("E26").Value = '''  +6.43%  '

$ws.Range("D27").Value = '''11.08'
$ws.Range("E27").Value = '''  +10.90%  '

$ws.Range("D28").Value = '''37.57'
$ws.Range("E28").Value = '''  +0.13%  '

$ws.Range("D29").Value = '''5.67'
$ws.Range("E29").Value = '''  -1.84%  '

$ws.Range("D30").Value = '''718.82'
$ws.Range("E30").Value = '''  +0.67%  '

$ws.Range("D31").Value = '''13.75'
$ws.Range("E31").Value = '''  +4.77%  '

$ws.Range("E32").Value = '''  +3.38%  '

$ws.Range("D33").Value = '''2.92'
$ws.Range("E33").Value = '''  +4.93%  '

$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").Value = '''0.0₃0890'
$ws.Range("E34").Value = '''  +5.15%  '

$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").Value = '''41.66'
$ws.Range("E35").Value = '''  +0.16%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '''6.09'
$ws.Range("E36").Value = '''  +15.23%  '

$ws.Range("D37").Value = '''60.64'
$ws.Range("E37").Value = '''  +4.21%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '''3.06'
$ws.Range("E38").Value = '''  +10.39%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = '''0.400'
$ws.Range("E39").Value = '''  +19.23%  '

$ws.Range("E40").Value = '''  -2.22%  '

$ws.Range("D41").Value = '''0.996'
$ws.Range("E41").Value = '''  -0.26%  '

$ws.Range("D42").Value = '''0.0497'
$ws.Range("E42").Value = '''  +7.74%  '

$ws.Range("D43").Value = '''3.11'
$ws.Range("E43").Value = '''  +3.15%  '

$ws.Range("D44").Value = '''2.97'
$ws.Range("E44").Value = '''  +1.88%  '

$ws.Range("E45").Value = '''  +3.40%  '

$ws.Range("D46").Value = '''3.36'
$ws.Range("E46").Value = '''  +6.08%  '

$ws.Range("E47").Value = '''  +0.24%  '

$ws.Range("E48").Value = '''  +1.57%  '

$ws.Range("E49").Value = '''  -0.17%  '

$ws.Range("D50").Value = '''145.12'
$ws.Range("E50").Value = '''  +0.43%  '

$ws.Range("E51").Value = '''  +29.88%  '

